$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the old 5-row sample block (rows 2-6) with a new 4-row data block (rows 2-5) ---

# Row 2
$ws.Cells.Item(2, 1).Value = 45040.50694444445
$ws.Cells.Item(2, 2).Value = 12.343
$ws.Cells.Item(2, 3).Value = 8.237
$ws.Cells.Item(2, 4).Value = 3.467
$ws.Cells.Item(2, 5).Value = 26.724
$ws.Cells.Item(2, 6).Value = 19.801
$ws.Cells.Item(2, 7).Value = 9.406000000000001
$ws.Cells.Item(2, 8).Value = 27.816
$ws.Cells.Item(2, 9).Value = 15.125
$ws.Cells.Item(2, 10).Value = 5.951
$ws.Cells.Item(2, 11).Value = 8.669
$ws.Cells.Item(2, 12).Value = 10.525
$ws.Cells.Item(2, 13).Value = 11.375
$ws.Cells.Item(2, 14).Value = 3.134
$ws.Cells.Item(2, 15).Value = 9.775
$ws.Cells.Item(2, 16).Value = 13.324
$ws.Cells.Item(2, 17).Value = 8.832000000000001
$ws.Cells.Item(2, 18).Value = 2.74
$ws.Cells.Item(2, 19).Value = 1.542
$ws.Cells.Item(2, 20).Value = 141.399
$ws.Cells.Item(2, 21).Value = 27.009
$ws.Cells.Item(2, 22).Value = 9.023
$ws.Cells.Item(2, 23).Value = 17.239
$ws.Cells.Item(2, 24).Value = 8.882
$ws.Cells.Item(2, 25).Value = 2.65
$ws.Cells.Item(2, 26).Value = 15.251
$ws.Cells.Item(2, 27).Value = 7.97
$ws.Cells.Item(2, 28).Value = 7.4
$ws.Cells.Item(2, 29).Value = 8.521000000000001
$ws.Cells.Item(2, 30).Value = 11.028
$ws.Cells.Item(2, 31).Value = 2.858
$ws.Cells.Item(2, 32).Value = 25.233
$ws.Cells.Item(2, 33).Value = 4.735
$ws.Cells.Item(2, 34).Value = 11.279

# Row 3
$ws.Cells.Item(3, 1).Value = 45040.51388888889
$ws.Cells.Item(3, 2).Value = 21.018
$ws.Cells.Item(3, 3).Value = 15.348
$ws.Cells.Item(3, 4).Value = 1.913
$ws.Cells.Item(3, 5).Value = 45.916
$ws.Cells.Item(3, 6).Value = 36.76
$ws.Cells.Item(3, 7).Value = 16.377
$ws.Cells.Item(3, 8).Value = 61.243
$ws.Cells.Item(3, 9).Value = 25.596
$ws.Cells.Item(3, 10).Value = 11.141
$ws.Cells.Item(3, 11).Value = 16.328
$ws.Cells.Item(3, 12).Value = 18.367
$ws.Cells.Item(3, 13).Value = 19.61
$ws.Cells.Item(3, 14).Value = 5.313
$ws.Cells.Item(3, 15).Value = 16.542
$ws.Cells.Item(3, 16).Value = 23.354
$ws.Cells.Item(3, 17).Value = 14.223
$ws.Cells.Item(3, 18).Value = 1.409
$ws.Cells.Item(3, 19).Value = 1.14
$ws.Cells.Item(3, 20).Value = 244.474
$ws.Cells.Item(3, 21).Value = 46.226
$ws.Cells.Item(3, 22).Value = 15.269
$ws.Cells.Item(3, 23).Value = 30.739
$ws.Cells.Item(3, 24).Value = 16.072
$ws.Cells.Item(3, 25).Value = 2.737
$ws.Cells.Item(3, 26).Value = 30.682
$ws.Cells.Item(3, 27).Value = 13.487
$ws.Cells.Item(3, 28).Value = 12.098
$ws.Cells.Item(3, 29).Value = 14.183
$ws.Cells.Item(3, 30).Value = 19.259
$ws.Cells.Item(3, 31).Value = 1.162
$ws.Cells.Item(3, 32).Value = 56.058
$ws.Cells.Item(3, 33).Value = 8.455
$ws.Cells.Item(3, 34).Value = 19.09

# Row 4
$ws.Cells.Item(4, 1).Value = 45040.52083333334
$ws.Cells.Item(4, 2).Value = 9.506
$ws.Cells.Item(4, 3).Value = 6.883
$ws.Cells.Item(4, 4).Value = 1.072
$ws.Cells.Item(4, 5).Value = 20.878
$ws.Cells.Item(4, 6).Value = 16.437
$ws.Cells.Item(4, 7).Value = 7.376
$ws.Cells.Item(4, 8).Value = 32.75
$ws.Cells.Item(4, 9).Value = 11.635
$ws.Cells.Item(4, 10).Value = 5.018
$ws.Cells.Item(4, 11).Value = 7.213
$ws.Cells.Item(4, 12).Value = 8.352
$ws.Cells.Item(4, 13).Value = 8.977
$ws.Cells.Item(4, 14).Value = 2.418
$ws.Cells.Item(4, 15).Value = 7.519
$ws.Cells.Item(4, 16).Value = 10.583
$ws.Cells.Item(4, 17).Value = 6.607
$ws.Cells.Item(4, 18).Value = 0.898
$ws.Cells.Item(4, 19).Value = 0.606
$ws.Cells.Item(4, 20).Value = 107.145
$ws.Cells.Item(4, 21).Value = 21.171
$ws.Cells.Item(4, 22).Value = 6.941
$ws.Cells.Item(4, 23).Value = 13.948
$ws.Cells.Item(4, 24).Value = 7.252
$ws.Cells.Item(4, 25).Value = 1.362
$ws.Cells.Item(4, 26).Value = 15.604
$ws.Cells.Item(4, 27).Value = 6.131
$ws.Cells.Item(4, 28).Value = 5.57
$ws.Cells.Item(4, 29).Value = 6.519
$ws.Cells.Item(4, 30).Value = 8.749000000000001
$ws.Cells.Item(4, 31).Value = 0.732
$ws.Cells.Item(4, 32).Value = 30.195
$ws.Cells.Item(4, 33).Value = 3.775
$ws.Cells.Item(4, 34).Value = 8.678000000000001

# Row 5
$ws.Cells.Item(5, 1).Value = 45040.52777777778
$ws.Cells.Item(5, 2).Value = 6.64
$ws.Cells.Item(5, 3).Value = 4.81
$ws.Cells.Item(5, 4).Value = 0.77
$ws.Cells.Item(5, 5).Value = 14.62
$ws.Cells.Item(5, 6).Value = 11.46
$ws.Cells.Item(5, 7).Value = 5.15
$ws.Cells.Item(5, 8).Value = 22.49
$ws.Cells.Item(5, 9).Value = 8.140000000000001
$ws.Cells.Item(5, 10).Value = 3.5
$ws.Cells.Item(5, 11).Value = 5.01
$ws.Cells.Item(5, 12).Value = 5.85
$ws.Cells.Item(5, 13).Value = 6.31
$ws.Cells.Item(5, 14).Value = 1.69
$ws.Cells.Item(5, 15).Value = 5.26
$ws.Cells.Item(5, 16).Value = 7.39
$ws.Cells.Item(5, 17).Value = 4.66
$ws.Cells.Item(5, 18).Value = 0.67
$ws.Cells.Item(5, 19).Value = 0.43
$ws.Cells.Item(5, 20).Value = 72.79000000000001
$ws.Cells.Item(5, 21).Value = 14.79
$ws.Cells.Item(5, 22).Value = 4.86
$ws.Cells.Item(5, 23).Value = 9.720000000000001
$ws.Cells.Item(5, 24).Value = 5.07
$ws.Cells.Item(5, 25).Value = 0.98
$ws.Cells.Item(5, 26).Value = 10.66
$ws.Cells.Item(5, 27).Value = 4.29
$ws.Cells.Item(5, 28).Value = 3.91
$ws.Cells.Item(5, 29).Value = 4.58
$ws.Cells.Item(5, 30).Value = 6.13
$ws.Cells.Item(5, 31).Value = 0.53
$ws.Cells.Item(5, 32).Value = 20.63
$ws.Cells.Item(5, 33).Value = 2.63
$ws.Cells.Item(5, 34).Value = 6.07

# Remove the now-unused trailing row (old data had 5 data rows, new data has 4)
$ws.Rows.Item(6).Delete()

# Column widths widen slightly to fit the new (longer) numeric values
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(5).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(9).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(20).ColumnWidth = 8.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17
